# Merge changes from remote repository
#
# Content change: column C ("Result") gets a value for the two data rows
# (row 2: testone@gmail.com, row 3: testtwo@gmail.com) recording the page
# title that was opened for the login link: "Facebook – log in or sign up".
# This introduces one new shared string and two new cells (C2, C3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$resultText = "Facebook " + [char]0x2013 + " log in or sign up"

$ws.Range("C2").Value = $resultText
$ws.Range("C3").Value = $resultText
